$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 539.8645649892381
$ws.Range("D2").Value = 121.3843902993817
$ws.Range("E2").Value = 386
$ws.Range("F2").Value = 447
$ws.Range("G2").Value = 505
$ws.Range("H2").Value = 595
$ws.Range("C3").Value = 43.71295319921021
$ws.Range("D3").Value = 4.81069944699828
$ws.Range("F3").Value = 40.42
$ws.Range("G3").Value = 43.65
$ws.Range("H3").Value = 46.93
$ws.Range("C4").Value = 1.634615819058292
$ws.Range("D4").Value = 2.899470628101848
$ws.Range("F4").Value = 0.59
$ws.Range("G4").Value = 1.24
$ws.Range("H4").Value = 2.22
$ws.Range("C5").Value = 323.2238041873454
$ws.Range("D5").Value = 10.29980121078036
$ws.Range("F5").Value = 317.83
$ws.Range("G5").Value = 325.28
$ws.Range("H5").Value = 331.07
$ws.Range("C6").Value = 22.25705242631233
$ws.Range("D6").Value = 1.91375967437697
$ws.Range("F6").Value = 21.12
$ws.Range("G6").Value = 22.05
$ws.Range("H6").Value = 22.95
$ws.Range("C7").Value = -76.67362852873686
$ws.Range("D7").Value = 23.63823620582192
$ws.Range("C8").Value = 7.572089589466981
$ws.Range("D8").Value = 6.818340091312328
$ws.Range("F8").Value = 7.8
$ws.Range("C9").Value = 9.323242969208602
$ws.Range("D9").Value = 1.688459102614802
$ws.Range("C10").Value = 867.8304938897485
$ws.Range("D10").Value = 0.4612111841487337
$ws.Range("C11").Value = 0.5571688973797961
$ws.Range("D11").Value = 0.5908796707351528
$ws.Range("C12").Value = 22.71885506163616
$ws.Range("D12").Value = 12.2909692008134
$ws.Range("G12").Value = 23
$ws.Range("H12").Value = 37
$ws.Range("I12").Value = 40
$ws.Range("C13").Value = 0.6726990945799314
$ws.Range("D13").Value = 0.7488468742878095
$ws.Range("C14").Value = 1.829617820232314
$ws.Range("D14").Value = 1.668499562659637
$ws.Range("C15").Value = 94.07362852873686
$ws.Range("D15").Value = 23.63823620581789
$ws.Range("C16").Value = -85.80559799117486
$ws.Range("D16").Value = 21.31908538193078
$ws.Range("F16").Value = -102.9618361134822
$ws.Range("H16").Value = -66.14699179957641
$ws.Range("C17").Value = -78.23350840170788
$ws.Range("D17").Value = 25.88644199732135
$ws.Range("F17").Value = -93.49305820175223
$ws.Range("G17").Value = -74.18978441047734
